# Generate Report for Handback
# Adds a newly-handed-back file (81be7453-d1f2-4473-bc59-0b9cc878ac7e.md) as a
# new row on all three sheets, and refreshes the existing row's data to
# reflect the latest handback run (new source file, new target xliff name,
# refreshed timestamps) for 7860255e-61b4-4612-a2bb-63b1b2ec886e.md.

$wb = $excel.ActiveWorkbook

$oldGuid = "3df035ea-ee68-4808-883d-efd6d55d3d47"
$guid1   = "7860255e-61b4-4612-a2bb-63b1b2ec886e"
$guid2   = "81be7453-d1f2-4473-bc59-0b9cc878ac7e"

$hash1 = "a0f5e76e9cc141511fa94a51d1f44f0335f4ca33"
$hash2 = "45b7a16c45c193b999653ed8f9ce75d7df741a0e"

$hyperlinkColor = 15570276  # RGB(100,149,237) == FF6495ED, matches the workbook's HyperLink style

function Apply-HyperlinkStyle($range) {
    $range.Font.Color = $hyperlinkColor
    $range.Font.Underline = $true
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

# Refresh row 2: the source file name for this row is now identified by the
# new GUID (the repo re-ran the handback against a renamed/regenerated file).
$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Range("B2").Value = "e2e\$guid1.md"

$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("B3").Value = "e2e\$guid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-31 03:10:34"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Apply-HyperlinkStyle($wsOverview.Range("B3"))
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/275565698c2f5a98d2bfb2122b0036a5dddde040/e2e/$guid2.md", $null, $null, "e2e\$guid2.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

# Refresh row 2 for the already-known file: new target xliff name + timestamps
$wsZh.Range("G2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-31 03:10:29"
$wsZh.Range("J2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 03:10:53"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# New row for the newly handed-back file
$rowZh = $loZh.ListRows.Add()

$wsZh.Range("A3").Value = "$guid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-31 03:10:29"
$wsZh.Range("I3").Value = "$guid2.md"
$wsZh.Range("J3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 03:10:53"
$wsZh.Range("M3").Value = "True"
$wsZh.Range("O3").Value = "False"

$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Apply-HyperlinkStyle($wsZh.Range("A3"))
Apply-HyperlinkStyle($wsZh.Range("I3"))

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/275565698c2f5a98d2bfb2122b0036a5dddde040/e2e/$guid2.md", $null, $null, "$guid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4c4d8dc38b0021116556a43c8a671bb964d336a1/e2e/$guid2.md", $null, $null, "$guid2.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

# Refresh row 2 for the already-known file: new target xliff name + timestamps
$wsDe.Range("G2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("J2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 03:11:01"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# New row for the newly handed-back file
$rowDe = $loDe.ListRows.Add()

$wsDe.Range("A3").Value = "$guid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-31 03:10:34"
$wsDe.Range("I3").Value = "$guid2.md"
$wsDe.Range("J3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 03:11:01"
$wsDe.Range("M3").Value = "True"
$wsDe.Range("O3").Value = "False"

$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Apply-HyperlinkStyle($wsDe.Range("A3"))
Apply-HyperlinkStyle($wsDe.Range("I3"))

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/275565698c2f5a98d2bfb2122b0036a5dddde040/e2e/$guid2.md", $null, $null, "$guid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1187d512b3dff162dc799705d09106d97509b8d7/e2e/$guid2.md", $null, $null, "$guid2.md") | Out-Null

Write-Output "Handback report generated."
